# 07_cues.xlsx -- lock in current version
# Replace the word/image/category cue triples (rows 2-49) with the refreshed
# stimulus set. Cells are written column-by-column (all of A, then all of B,
# then all of C) in row order, matching the order in which the authoring tool
# originally interned these shared strings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A -- word
$ws.Cells.Item(2, 1).Value = 'wenden'
$ws.Cells.Item(3, 1).Value = 'antun'
$ws.Cells.Item(4, 1).Value = 'streichen'
$ws.Cells.Item(5, 1).Value = 'rufen'
$ws.Cells.Item(6, 1).Value = 'segeln'
$ws.Cells.Item(7, 1).Value = 'holen'
$ws.Cells.Item(8, 1).Value = 'enden'
$ws.Cells.Item(9, 1).Value = 'steuern'
$ws.Cells.Item(10, 1).Value = 'zahlen'
$ws.Cells.Item(11, 1).Value = 'opfern'
$ws.Cells.Item(12, 1).Value = 'tanzen'
$ws.Cells.Item(13, 1).Value = 'küssen'
$ws.Cells.Item(14, 1).Value = 'schalten'
$ws.Cells.Item(15, 1).Value = 'altern'
$ws.Cells.Item(16, 1).Value = 'grüßen'
$ws.Cells.Item(17, 1).Value = 'klagen'
$ws.Cells.Item(18, 1).Value = 'sammeln'
$ws.Cells.Item(19, 1).Value = 'schaden'
$ws.Cells.Item(20, 1).Value = 'sparen'
$ws.Cells.Item(21, 1).Value = 'beißen'
$ws.Cells.Item(22, 1).Value = 'bitten'
$ws.Cells.Item(23, 1).Value = 'hören'
$ws.Cells.Item(24, 1).Value = 'kehren'
$ws.Cells.Item(25, 1).Value = 'husten'
$ws.Cells.Item(26, 1).Value = 'drohen'
$ws.Cells.Item(27, 1).Value = 'dürfen'
$ws.Cells.Item(28, 1).Value = 'lassen'
$ws.Cells.Item(29, 1).Value = 'orten'
$ws.Cells.Item(30, 1).Value = 'dringen'
$ws.Cells.Item(31, 1).Value = 'buchen'
$ws.Cells.Item(32, 1).Value = 'weigern'
$ws.Cells.Item(33, 1).Value = 'platzen'
$ws.Cells.Item(34, 1).Value = 'spenden'
$ws.Cells.Item(35, 1).Value = 'ächzen'
$ws.Cells.Item(36, 1).Value = 'treiben'
$ws.Cells.Item(37, 1).Value = 'wachsen'
$ws.Cells.Item(38, 1).Value = 'kosten'
$ws.Cells.Item(39, 1).Value = 'süßen'
$ws.Cells.Item(40, 1).Value = 'jubeln'
$ws.Cells.Item(41, 1).Value = 'dauern'
$ws.Cells.Item(42, 1).Value = 'heilen'
$ws.Cells.Item(43, 1).Value = 'quälen'
$ws.Cells.Item(44, 1).Value = 'stören'
$ws.Cells.Item(45, 1).Value = 'pflegen'
$ws.Cells.Item(46, 1).Value = 'füttern'
$ws.Cells.Item(47, 1).Value = 'bremsen'
$ws.Cells.Item(48, 1).Value = 'stärken'
$ws.Cells.Item(49, 1).Value = 'drücken'

# Column B -- image
$ws.Cells.Item(2, 2).Value = 'none'
$ws.Cells.Item(3, 2).Value = 'face/face023.jpg'
$ws.Cells.Item(4, 2).Value = 'house/house010.jpg'
$ws.Cells.Item(5, 2).Value = 'none'
$ws.Cells.Item(6, 2).Value = 'face/face022.jpg'
$ws.Cells.Item(7, 2).Value = 'house/house013.jpg'
$ws.Cells.Item(8, 2).Value = 'none'
$ws.Cells.Item(9, 2).Value = 'face/face019.jpg'
$ws.Cells.Item(10, 2).Value = 'house/house002.jpg'
$ws.Cells.Item(11, 2).Value = 'none'
$ws.Cells.Item(12, 2).Value = 'house/house014.jpg'
$ws.Cells.Item(13, 2).Value = 'face/face002.jpg'
$ws.Cells.Item(14, 2).Value = 'none'
$ws.Cells.Item(15, 2).Value = 'face/face027.jpg'
$ws.Cells.Item(16, 2).Value = 'house/house020.jpg'
$ws.Cells.Item(17, 2).Value = 'none'
$ws.Cells.Item(18, 2).Value = 'face/face025.jpg'
$ws.Cells.Item(19, 2).Value = 'face/face018.jpg'
$ws.Cells.Item(20, 2).Value = 'none'
$ws.Cells.Item(21, 2).Value = 'house/house023.jpg'
$ws.Cells.Item(22, 2).Value = 'face/face005.jpg'
$ws.Cells.Item(23, 2).Value = 'none'
$ws.Cells.Item(24, 2).Value = 'face/face014.jpg'
$ws.Cells.Item(25, 2).Value = 'house/house030.jpg'
$ws.Cells.Item(26, 2).Value = 'none'
$ws.Cells.Item(27, 2).Value = 'house/house024.jpg'
$ws.Cells.Item(28, 2).Value = 'face/face006.jpg'
$ws.Cells.Item(29, 2).Value = 'none'
$ws.Cells.Item(30, 2).Value = 'house/house015.jpg'
$ws.Cells.Item(31, 2).Value = 'house/house009.jpg'
$ws.Cells.Item(32, 2).Value = 'none'
$ws.Cells.Item(33, 2).Value = 'face/face026.jpg'
$ws.Cells.Item(34, 2).Value = 'house/house011.jpg'
$ws.Cells.Item(35, 2).Value = 'none'
$ws.Cells.Item(36, 2).Value = 'face/face001.jpg'
$ws.Cells.Item(37, 2).Value = 'face/face013.jpg'
$ws.Cells.Item(38, 2).Value = 'none'
$ws.Cells.Item(39, 2).Value = 'face/face030.jpg'
$ws.Cells.Item(40, 2).Value = 'face/face009.jpg'
$ws.Cells.Item(41, 2).Value = 'none'
$ws.Cells.Item(42, 2).Value = 'house/house028.jpg'
$ws.Cells.Item(43, 2).Value = 'house/house016.jpg'
$ws.Cells.Item(44, 2).Value = 'none'
$ws.Cells.Item(45, 2).Value = 'house/house029.jpg'
$ws.Cells.Item(46, 2).Value = 'house/house025.jpg'
$ws.Cells.Item(47, 2).Value = 'none'
$ws.Cells.Item(48, 2).Value = 'face/face010.jpg'
$ws.Cells.Item(49, 2).Value = 'house/house012.jpg'

# Column C -- category
$ws.Cells.Item(2, 3).Value = 'none'
$ws.Cells.Item(3, 3).Value = 'face'
$ws.Cells.Item(4, 3).Value = 'house'
$ws.Cells.Item(5, 3).Value = 'none'
$ws.Cells.Item(6, 3).Value = 'face'
$ws.Cells.Item(7, 3).Value = 'house'
$ws.Cells.Item(8, 3).Value = 'none'
$ws.Cells.Item(9, 3).Value = 'face'
$ws.Cells.Item(10, 3).Value = 'house'
$ws.Cells.Item(11, 3).Value = 'none'
$ws.Cells.Item(12, 3).Value = 'house'
$ws.Cells.Item(13, 3).Value = 'face'
$ws.Cells.Item(14, 3).Value = 'none'
$ws.Cells.Item(15, 3).Value = 'face'
$ws.Cells.Item(16, 3).Value = 'house'
$ws.Cells.Item(17, 3).Value = 'none'
$ws.Cells.Item(18, 3).Value = 'face'
$ws.Cells.Item(19, 3).Value = 'face'
$ws.Cells.Item(20, 3).Value = 'none'
$ws.Cells.Item(21, 3).Value = 'house'
$ws.Cells.Item(22, 3).Value = 'face'
$ws.Cells.Item(23, 3).Value = 'none'
$ws.Cells.Item(24, 3).Value = 'face'
$ws.Cells.Item(25, 3).Value = 'house'
$ws.Cells.Item(26, 3).Value = 'none'
$ws.Cells.Item(27, 3).Value = 'house'
$ws.Cells.Item(28, 3).Value = 'face'
$ws.Cells.Item(29, 3).Value = 'none'
$ws.Cells.Item(30, 3).Value = 'house'
$ws.Cells.Item(31, 3).Value = 'house'
$ws.Cells.Item(32, 3).Value = 'none'
$ws.Cells.Item(33, 3).Value = 'face'
$ws.Cells.Item(34, 3).Value = 'house'
$ws.Cells.Item(35, 3).Value = 'none'
$ws.Cells.Item(36, 3).Value = 'face'
$ws.Cells.Item(37, 3).Value = 'face'
$ws.Cells.Item(38, 3).Value = 'none'
$ws.Cells.Item(39, 3).Value = 'face'
$ws.Cells.Item(40, 3).Value = 'face'
$ws.Cells.Item(41, 3).Value = 'none'
$ws.Cells.Item(42, 3).Value = 'house'
$ws.Cells.Item(43, 3).Value = 'house'
$ws.Cells.Item(44, 3).Value = 'none'
$ws.Cells.Item(45, 3).Value = 'house'
$ws.Cells.Item(46, 3).Value = 'house'
$ws.Cells.Item(47, 3).Value = 'none'
$ws.Cells.Item(48, 3).Value = 'face'
$ws.Cells.Item(49, 3).Value = 'house'

